# Updated cryptos list on Thu Oct 12 09:50:36 UTC 2023 with GitHub Actions
# Refreshes the Price / Volume(1h) columns (and, where the rank order
# changed, the Coin name + Link) for the scraped cryptocurrency table.
#
# Every text value is written with a leading apostrophe ("'") so Excel
# stores it as literal text (matching the original inlineStr cells,
# which include things like "26.774.02" and " -1.59% " that would
# otherwise be re-interpreted as numbers). The Style reset back to
# "Normal" immediately afterwards clears the quote-prefix flag that
# the leading apostrophe leaves behind, so the cell's formatting stays
# identical to before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.774.02"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  -1.59%  "
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(3, 4).Value = "'1.551.09"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  -1.45%  "
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "'  +0.22%  "
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(5, 4).Value = "'204.57"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  -1.23%  "
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(6, 4).Value = "'0.481"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  -1.61%  "
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'  +0.29%  "
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(8, 2).Value = "'Solana"
$ws.Cells.Item(8, 2).Style = "Normal"
$ws.Cells.Item(8, 3).Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(8, 3).Style = "Normal"
$ws.Cells.Item(8, 4).Value = "'21.43"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "'  -4.15%  "
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(9, 2).Value = "'Cardano"
$ws.Cells.Item(9, 2).Style = "Normal"
$ws.Cells.Item(9, 3).Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(9, 3).Style = "Normal"
$ws.Cells.Item(9, 4).Value = "'0.245"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'  -0.95%  "
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  -1.47%  "
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'  -0.70%  "
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(12, 4).Value = "'1.773.72"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'  -1.29%  "
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(13, 4).Value = "'1.558.52"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  -0.99%  "
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(14, 4).Value = "'3.67"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  -2.56%  "
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  -1.95%  "
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(16, 4).Value = "'26.769.41"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  -1.56%  "
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(17, 4).Value = "'61.12"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  -2.38%  "
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(18, 4).Value = "'213.68"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  -0.75%  "
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  -0.42%  "
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(20, 4).Value = "'0.0₃0680"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  -0.93%  "
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  +0.20%  "
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  -1.62%  "
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(23, 4).Value = "'9.07"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  -3.66%  "
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(24, 4).Value = "'2.00"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  -0.11%  "
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(25, 4).Value = "'151.92"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  -0.54%  "
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'  -1.94%  "
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(27, 4).Value = "'14.86"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  -0.63%  "
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  +0.19%  "
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  -2.32%  "
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  -0.31%  "
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  -3.04%  "
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(32, 4).Value = "'3.16"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  -0.60%  "
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Cells.Item(33, 4).Value = "'1.364.37"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "'  -2.60%  "
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  -0.84%  "
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Cells.Item(35, 4).Value = "'1.49"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  -4.49%  "
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  -0.40%  "
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(37, 4).Value = "'0.919"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  -2.78%  "
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  -2.26%  "
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(39, 4).Value = "'0.520"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  +0.07%  "
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(40, 4).Value = "'0.802"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  -2.01%  "
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'  +0.21%  "
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(42, 2).Value = "'WEMIXToken"
$ws.Cells.Item(42, 2).Style = "Normal"
$ws.Cells.Item(42, 3).Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(42, 3).Style = "Normal"
$ws.Cells.Item(42, 4).Value = "'0.988"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  -0.16%  "
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(43, 2).Value = "'FraxShare"
$ws.Cells.Item(43, 2).Style = "Normal"
$ws.Cells.Item(43, 3).Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(43, 3).Style = "Normal"
$ws.Cells.Item(43, 4).Value = "'5.54"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  +3.54%  "
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  +0.11%  "
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  -2.89%  "
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(46, 4).Value = "'62.93"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'  -1.40%  "
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  -2.55%  "
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(48, 4).Value = "'1.687.63"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  -1.20%  "
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(49, 4).Value = "'86.10"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'  +0.13%  "
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(50, 4).Value = "'0.0511"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  +3.58%  "
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(51, 4).Value = "'0.0₇0978"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'  -1.09%  "
$ws.Cells.Item(51, 5).Style = "Normal"
